# Update "想去人数" (number interested) counts for several events.
# Values increase by a small amount in both the "展览" sheet (rows 3,10,16,22,29,30,31,33)
# and the "全部类型" sheet (same events, but row 33 is shifted to row 34 there
# because "全部类型" also includes the single row from "演出").

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 3;  Value = 564 },
    @{ Row = 10; Value = 16457 },
    @{ Row = 16; Value = 128 },
    @{ Row = 22; Value = 38 },
    @{ Row = 29; Value = 900 },
    @{ Row = 30; Value = 59 },
    @{ Row = 31; Value = 5063 }
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}
$ws1.Cells.Item(33, 6).Value = 11375

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
$ws4.Cells.Item(34, 6).Value = 11375
